# #5: property boat&car done
#
# The "汽車" (car) sheet only had a 6-column stub (name/year/owner/
# register_reason/register_date placeholders). Bring it in line with the
# other two property sheets (土地/建物) which use the full 14-column
# layout (name, capacity, owner, register_date, register_reason,
# acquire_value, property_category, category, date, legislator_name,
# legislator_id, source_file, index) and fill in the actual car record
# (VOLKSWAGEN Passat, 1984cc) plus its 1740cc engine/legislator id quirk
# that mirrors the other sheets' data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row (row 1) --------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Carry the existing header style (bold font + border, centered) from B1
# onto the newly-added header cells C1:N1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data row (row 2) -----------------------------------------------------
$ws.Range("A2").Value = 33
$ws.Range("B2").Value = "VOLKSWAGENPSDDAT2.0"
$ws.Range("C2").Value = 1984
$ws.Range("D2").Value = "林世嘉"
$ws.Range("E2").Value = "94年01月28曰"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = "(超過五年）"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# "2012-05-01" looks like a date literal, and a plain .Value assignment
# would get auto-coerced into a date serial number. Build it as a text
# formula result in a scratch cell first, then paste just the computed
# value back in so the cell ends up holding the literal text string.
$ws.Range("ZZ1").Formula = '="2012-05-01"'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4163) | Out-Null
$ws.Range("ZZ1").Clear() | Out-Null
$excel.CutCopyMode = $false

$ws.Range("K2").Value = "林世嘉"
$ws.Range("L2").Value = 1740
$ws.Range("M2").Value = "tmpada11"
$ws.Range("N2").Value = 33

# Carry the existing data-row style from B2 onto the newly-added data
# cells C2:N2 (this also re-applies the plain/no-border style to J2,
# which the value-paste step above left without formatting).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2:N2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
